# "Separate tichonov method to another file"
# Column A (rows 2..31) switches from literal values to a formula that
# recomputes the same numbers from columns D/E/F:
#   A{r} = (4 * D{r}) + (E{r} * F{r})
# A3:A31 is entered as one fill so the engine records it as a shared formula
# (matching the t="shared" ref="A3:A31" si="0" pattern in the target file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Formula = "=(4*D2)+(E2*F2)"
$ws.Range("A3:A31").Formula = "=(4*D3)+(E3*F3)"

# Selection moved from J10 to M30 before the file was last saved.
$ws.Range("M30").Select()
